$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "maa://25390 (95.35), maa://24702 (94.77), maa://36681 (86.42)"
$ws.Range("L2").Value = "maa://39402 (94.23), *maa://24633 (56.1), *maa://30515 (70.09), *maa://34787 (73.56), ***maa://20792 (11.93), ***maa://29083 (27.78)"
$ws.Range("T2").Value = "maa://22742 (91.33), *maa://20791 (62.2)"
$ws.Range("AB2").Value = "maa://21246 (91.42), maa://36684 (93.29), ***maa://22731 (6.25)"
$ws.Range("AF2").Value = "maa://25251 (91.79), **maa://21730 (30.49), ***maa://39501 (11.36), **maa://36675 (50.0)"
$ws.Range("H3").Value = "maa://21247 (98.37), *maa://22748 (60.0)"
$ws.Range("L3").Value = "*maa://22880 (65.02), maa://20276 (88.0), *maa://22749 (75.0)"
$ws.Range("P3").Value = "maa://21249 (94.51), maa://26254 (97.14), **maa://22738 (50.0)"
$ws.Range("T3").Value = "maa://24617 (90.55), maa://45854 (85.29), **maa://20790 (43.48), ***maa://37170 (16.67)"
$ws.Range("X3").Value = "maa://27396 (84.35), maa://27484 (96.21), maa://27480 (83.78)"
$ws.Range("P4").Value = "maa://49983 (94.34), maa://50121 (92.0)"
$ws.Range("T4").Value = "maa://32509 (94.53), maa://27295 (88.24), maa://22754 (89.19), *maa://31008 (79.07), *maa://21746 (54.55)"
$ws.Range("X4").Value = "**maa://32495 (48.4), maa://43217 (92.92), ***maa://31785 (22.22), ***maa://36683 (29.79)"
$ws.Range("AF4").Value = "*maa://30062 (65.45), *maa://39394 (65.52), ***maa://26209 (13.04)"
$ws.Range("D5").Value = "maa://21245 (84.93), maa://22744 (81.48)"
$ws.Range("D6").Value = "maa://42407 (94.81)"
$ws.Range("L6").Value = "maa://24839 (98.8)"
$ws.Range("AF6").Value = "*maa://33152 (63.79), ***maa://22770 (26.09)"
$ws.Range("L7").Value = "maa://28624 (93.53), maa://24957 (97.78)"
$ws.Range("AF7").Value = "*maa://26191 (67.82), *maa://36671 (69.23), maa://45272 (96.67), *maa://42530 (62.5)"
$ws.Range("A8").Value = "更新日期：2025.05.11 13:19:54"
$ws.Range("P8").Value = "maa://32931 (84.09), *maa://21916 (62.69), maa://23252 (91.55), maa://37496 (97.22), **maa://22759 (45.45)"
$ws.Range("P9").Value = "maa://22736 (82.88)"
$ws.Range("T9").Value = "**maa://22866 (30.19), maa://26222 (98.28)"
$ws.Range("AB9").Value = "maa://28711 (86.76), **maa://39938 (45.0), **maa://45044 (40.0), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (97.06)"
$ws.Range("L10").Value = "**maa://24395 (38.89)"
$ws.Range("T10").Value = "maa://27395 (96.52), maa://22755 (88.62), **maa://22756 (40.91), ***maa://21737 (11.76)"
$ws.Range("X10").Value = "maa://22301 (97.85), maa://45828 (92.68), maa://22726 (100.0)"
$ws.Range("X12").Value = "maa://22753 (91.67), *maa://21485 (75.51), maa://37962 (91.67)"
$ws.Range("AB12").Value = "maa://23669 (95.64), maa://36677 (95.12), maa://39872 (92.86)"
$ws.Range("AF12").Value = "*maa://28932 (79.53), *maa://20106 (64.29), *maa://22769 (64.29)"
$ws.Range("D13").Value = "maa://24999 (92.57), maa://36673 (92.94), maa://25001 (86.3)"
$ws.Range("P13").Value = "maa://22676 (93.53), *maa://22583 (77.63), *maa://22500 (59.57), maa://48321 (100.0)"
$ws.Range("AF13").Value = "**maa://22737 (38.22), maa://39883 (87.72), *maa://39885 (52.63)"
$ws.Range("L14").Value = "maa://39841 (94.59), maa://26245 (96.53), maa://21288 (96.3), maa://36682 (95.56)"
$ws.Range("P14").Value = "maa://23250 (98.81), maa://20107 (87.1), maa://22772 (100.0), *maa://22745 (66.67)"
$ws.Range("AB14").Value = "maa://22764 (96.25)"
$ws.Range("D15").Value = "*maa://22743 (79.06), maa://22734 (84.43), *maa://30808 (65.22), *maa://36048 (64.58), maa://45058 (84.21)"
$ws.Range("H15").Value = "maa://24304 (86.67), maa://21478 (90.0)"
$ws.Range("AF15").Value = "maa://21364 (81.27), maa://36666 (81.12), *maa://22766 (68.03)"
$ws.Range("AF16").Value = "*maa://23911 (68.1), maa://27755 (93.88)"
$ws.Range("H17").Value = "maa://22430 (89.42), maa://39599 (85.14)"
$ws.Range("AF17").Value = "maa://50136 (94.44)"
$ws.Range("H18").Value = "maa://24421 (87.82)"
$ws.Range("L18").Value = "maa://22466 (92.08), *maa://22732 (51.38), maa://52226 (100.0)"
$ws.Range("D20").Value = "maa://21432 (90.82), maa://25198 (94.4), *maa://20795 (50.76), maa://36680 (91.18)"
$ws.Range("H20").Value = "maa://22864 (90.5), ***maa://43283 (25.0), maa://53361 (100.0)"
$ws.Range("L20").Value = "maa://41331 (85.71)"
$ws.Range("X20").Value = "maa://50085 (88.64), maa://49976 (85.71)"
$ws.Range("X21").Value = "maa://20110 (86.76), maa://34946 (93.1)"
$ws.Range("AB21").Value = "maa://21443 (82.56), ***maa://23820 (30.0), **maa://52223 (38.46)"
$ws.Range("AF21").Value = "maa://22524 (90.95), maa://22432 (81.98)"
$ws.Range("L22").Value = "*maa://27127 (78.46), *maa://22751 (71.23)"
$ws.Range("X22").Value = "maa://21282 (98.7), *maa://37649 (69.44)"
$ws.Range("L23").Value = "maa://39756 (95.41), maa://39875 (94.67)"
$ws.Range("X24").Value = "maa://29988 (84.1), maa://23504 (93.53), **maa://22892 (41.29), *maa://25141 (77.37), *maa://36663 (77.78), ***maa://22815 (23.08), maa://52227 (100.0)"
$ws.Range("AF24").Value = "maa://22523 (82.71), *maa://36672 (76.19), maa://29910 (93.75), **maa://21440 (35.71), maa://45831 (84.62)"
$ws.Range("D25").Value = "maa://29753 (95.3)"
$ws.Range("AB25").Value = "maa://31215 (89.47), maa://24516 (80.22), maa://26001 (84.48)"
$ws.Range("AF25").Value = "maa://20108 (95.8), maa://24621 (97.14), maa://36676 (97.06), maa://22771 (85.71), *maa://37772 (66.67)"
$ws.Range("H26").Value = "maa://24913 (91.75)"
$ws.Range("AB26").Value = "maa://42235 (95.62)"
$ws.Range("X28").Value = "maa://39929 (91.79), maa://41749 (91.67), ***maa://39723 (13.89)"
$ws.Range("AF28").Value = "maa://36660 (92.65), *maa://36701 (66.67)"
$ws.Range("H29").Value = "*maa://25175 (60.66)"
$ws.Range("L29").Value = "maa://28432 (93.92), maa://28440 (83.46), maa://31400 (98.88), *maa://28650 (71.43)"
$ws.Range("AF29").Value = "*maa://24080 (69.17), maa://42865 (81.0), ***maa://34960 (8.33)"
$ws.Range("D30").Value = "maa://45792 (88.89)"
$ws.Range("P30").Value = "maa://21442 (98.76)"
$ws.Range("AB30").Value = "maa://42979 (96.76), maa://45822 (100.0), maa://45045 (83.33)"
$ws.Range("L31").Value = "maa://35926 (93.85), maa://36258 (86.75), *maa://43904 (71.43)"
$ws.Range("H32").Value = "maa://21895 (97.51), maa://36667 (98.0), **maa://20793 (38.0), maa://22760 (100.0)"
$ws.Range("P34").Value = "maa://48817 (96.08)"
$ws.Range("AF38").Value = "maa://36697 (88.77)"
$ws.Range("H39").Value = "maa://36670 (88.89), maa://25199 (85.22), maa://30434 (92.56), maa://45059 (80.65), ***maa://25036 (18.52), *maa://44165 (66.67)"
$ws.Range("T39").Value = "maa://47079 (94.87), *maa://45788 (76.32), *maa://45790 (79.17)"
$ws.Range("P41").Value = "**maa://35616 (40.0), maa://43177 (93.33)"
$ws.Range("H43").Value = "maa://22525 (88.39), maa://21284 (87.04)"
$ws.Range("H44").Value = "maa://29768 (98.22), maa://27728 (96.15)"
$ws.Range("H46").Value = "maa://35931 (92.48), maa://43901 (94.0)"
$ws.Range("H47").Value = "maa://27410 (96.82), maa://29661 (97.52), maa://28038 (84.62)"
$ws.Range("H53").Value = "maa://32534 (94.78), **maa://32434 (33.33)"
$ws.Range("H55").Value = "maa://32532 (92.92)"
$ws.Range("H62").Value = "maa://42981 (96.3), maa://43903 (100.0)"

# G20 target value "3" must remain a TEXT string (matches source data semantics) while
# keeping the original cell style (numFmtId General). Directly assigning Range.Value = "3"
# would auto-coerce to a number (standard Excel behavior). Copying from C2, which already
# holds the literal text "3" with the identical style, preserves both the string type and style.
$ws.Range("C2").Copy($ws.Range("G20"))
